$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.27868390083313
$ws.Range("B1").Value = 2.371200084686279
$ws.Range("C1").Value = 2.466561079025269
$ws.Range("D1").Value = 3.246614694595337
$ws.Range("E1").Value = 2.262989521026611
